$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-7, columns D (Fecha), M (Volumen), N (Precio minimo),
# O (Precio maximo), P (Precio promedio ponderado), S (Precio $/Kg).
# This reflects a cyclic rotation of the weekly data rows (new row r <- old row r+4 mod 6).
$data = @{
    2 = @{ D = 44216; M = 200; N = 3500;  O = 4000;  P = 3750;  S = 1875 }
    3 = @{ D = 44216; M = 100; N = 3000;  O = 3000;  P = 3000;  S = 1500 }
    4 = @{ D = 44532; M = 100; N = 10000; O = 10000; P = 10000; S = 5000 }
    5 = @{ D = 44532; M = 100; N = 8000;  O = 8000;  P = 8000;  S = 4000 }
    6 = @{ D = 44195; M = 200; N = 3000;  O = 3500;  P = 3250;  S = 1625 }
    7 = @{ D = 44195; M = 100; N = 2500;  O = 2500;  P = 2500;  S = 1250 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("N$row").Value = $vals.N
    $ws.Range("O$row").Value = $vals.O
    $ws.Range("P$row").Value = $vals.P
    $ws.Range("S$row").Value = $vals.S
}

$wb.Save()
